$wb = $excel.ActiveWorkbook

# --- Sheet "scenarios" (sheet1) ---
$ws1 = $wb.Worksheets.Item("scenarios")

# Update row 2 values (B2, C2) with the new Linz AG scenario data
$ws1.Range("B2").Value = "Portfolio Status Quo"
$ws1.Range("C2").Value = "portfolio_quo"

# Remove the old row 3 (Szenario No Nuc Low C02 / scenario_lzn) entirely
$ws1.Rows.Item(3).Delete()

# --- Sheet "sub-scenarios" (sheet2) ---
$ws2 = $wb.Worksheets.Item("sub-scenarios")

# Header: B1 now holds "Prefix" instead of "Low Price Scenario"
$ws2.Range("B1").Value = "Prefix"

# Update the scenario names in column C
$ws2.Range("C2").Value = "Szenario Quo"
$ws2.Range("C3").Value = "Szenario Niedrig"
$ws2.Range("C4").Value = "Szenario Moderat"
$ws2.Range("C5").Value = "Scenario Hoch"

# --- Restore the selections recorded in each sheet view ---
$ws2.Activate()
$ws2.Range("G19").Select()

$ws1.Activate()
$ws1.Range("D11").Select()
